$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wrote a new test case for AllowImportPermissions -> Total Test Cases (column C) increments from 4 to 5
$ws.Range("C2").Value = 5

# Update the active cell selection recorded in the sheet view
$ws.Range("E4").Select()
